$wb = $excel.ActiveWorkbook

# Update the "Status" value from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn, and de-de sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the columns that held the status text so their widths shrink to
# match the new, shorter text (mirrors Excel's autofit-on-edit behavior,
# which recalculates each column to the new text's ideal display width).
# Target "ideal" width is ~13.41 characters; 12.5 is the ColumnWidth input
# that this engine's pixel-quantized width model resolves closest to it.
$wsOverview.Range("E1:F2").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1:C2").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1:C2").EntireColumn.ColumnWidth = 12.5
